# Apply the "Add idle anim for bat" workbook edit (Equipments sheet: new
# Armor rows + unhide previously-filtered rows + freeze header row).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Equipments")

# --- New "Armor" equipment rows (24-45) -----------------------------------
# Row 24: first Armor entry, fully populated (Name + Modifier text too).
$ws.Cells.Item(24, 1).Value = "Armor"
$ws.Cells.Item(24, 2).Value = "Common"
$ws.Cells.Item(24, 3).Value = "Leather Armor"
$ws.Cells.Item(24, 4).Value = "+ 10 Damage"

# Row 25: second Armor entry, also fully populated.
$ws.Cells.Item(25, 1).Value = "Armor"
$ws.Cells.Item(25, 2).Value = "Common"
$ws.Cells.Item(25, 3).Value = "Steal Armor"
$ws.Cells.Item(25, 4).Value = "+ 50 HP"

# Rows 26-30: Common rarity placeholders (Category + Rarity only).
26..30 | ForEach-Object {
    $ws.Cells.Item($_, 1).Value = "Armor"
    $ws.Cells.Item($_, 2).Value = "Common"
}

# Rows 31-38: Uncommon rarity placeholders.
31..38 | ForEach-Object {
    $ws.Cells.Item($_, 1).Value = "Armor"
    $ws.Cells.Item($_, 2).Value = "Uncommon"
}

# Rows 39-42: Rare rarity placeholders.
39..42 | ForEach-Object {
    $ws.Cells.Item($_, 1).Value = "Armor"
    $ws.Cells.Item($_, 2).Value = "Rare"
}

# Rows 43-44: Epic rarity placeholders.
43..44 | ForEach-Object {
    $ws.Cells.Item($_, 1).Value = "Armor"
    $ws.Cells.Item($_, 2).Value = "Epic"
}

# Row 45: Legendary rarity placeholder.
$ws.Cells.Item(45, 1).Value = "Armor"
$ws.Cells.Item(45, 2).Value = "Legendary"

# --- Unhide the rows that were previously hidden by the autofilter --------
$ws.Range("A2:A23").EntireRow.Hidden = $false

# --- Clear the column-E "blank" filter criterion, keep the filter buttons -
$ws.ShowAllData() | Out-Null

# --- Freeze the header row and set the active selection -------------------
$ws.Activate() | Out-Null
$ws.Range("A2").Select() | Out-Null
$excel.ActiveWindow.FreezePanes = $true
$ws.Range("C26").Select() | Out-Null
